$wb = $excel.ActiveWorkbook

# --- Sheet "Mensal" (sheet1) updates ---
$wsMensal = $wb.Worksheets.Item("Mensal")

# Row 13: update cota_obs_cm and anomalia
$wsMensal.Range("B13").Value = 30.72
$wsMensal.Range("D13").Value = -78.22

# Row 14: update date, cota_obs_cm and anomalia
$wsMensal.Range("A14").Value = 44066
$wsMensal.Range("B14").Value = 71.16
$wsMensal.Range("D14").Value = -37.54

# --- Sheet "Diario" (sheet2) updates ---
$wsDiario = $wb.Worksheets.Item("Diario")

# Row 365: update cota_obs_cm and anomalia
$wsDiario.Range("B365").Value = 21.96
$wsDiario.Range("D365").Value = -84.43000000000001

# Row 367: update cota_obs_cm and anomalia
$wsDiario.Range("B367").Value = 45.76
$wsDiario.Range("D367").Value = -67.55

# New rows 384-390: copy formatting (date style on col A, borders, etc.) from the
# last existing row (383), then overwrite with the new values below.
$wsDiario.Range("A383:D383").Copy()
$wsDiario.Range("A384:D390").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$newRows = @(
    @(44060, 131.07, 113.94, 15.04),
    @(44061, 120.43, 113.94, 5.7),
    @(44062, 141.56, 113.94, 24.25),
    @(44063, 182.5, 113.94, 60.18),
    @(44064, 195.12, 113.94, 71.26000000000001),
    @(44065, 178.66, 113.94, 56.8),
    @(44066, 156.07, 113.94, 36.98)
)

$r = 384
foreach ($row in $newRows) {
    $wsDiario.Range("A$r").Value = $row[0]
    $wsDiario.Range("B$r").Value = $row[1]
    $wsDiario.Range("C$r").Value = $row[2]
    $wsDiario.Range("D$r").Value = $row[3]
    $r = $r + 1
}
